# Apply "updates to growth modes and rates" changes to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New cell E2 = "model"
$ws.Range("E2").Value = "model"

# C3: "experiment" -> "model"
$ws.Range("C3").Value = "model"

# New cell E3 = "model"
$ws.Range("E3").Value = "model"

# New cell E4 = "experiment"
$ws.Range("E4").Value = "experiment"

# Update the active cell selection to D4, matching the saved view state
$ws.Activate()
$ws.Range("D4").Select()
